$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.339.04"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "'3.131.36"
$ws.Range("E3").Value = "  -4.91%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'211.63"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("D6").Value = "'622.94"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").Value = "'0.394"
$ws.Range("E7").Value = "  -6.37%  "
$ws.Range("D8").Value = "'0.715"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'3.132.96"
$ws.Range("E10").Value = "  -4.96%  "
$ws.Range("D11").Value = "'0.549"
$ws.Range("E11").Value = "  -8.45%  "
$ws.Range("D12").Value = "'0.180"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  -8.84%  "
$ws.Range("D14").Value = "'89.689.30"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'3.734.86"
$ws.Range("E15").Value = "  -3.90%  "
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").Value = "'5.23"
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("D17").Value = "'31.71"
$ws.Range("E17").Value = "  -7.90%  "
$ws.Range("D18").Value = "'3.212.34"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'3.29"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").Value = "'0.0000210"
$ws.Range("E20").Value = "  +12.09%  "
$ws.Range("D21").Value = "'13.08"
$ws.Range("E21").Value = "  -8.25%  "
$ws.Range("D22").Value = "'422.35"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").Value = "'8.29"
$ws.Range("E23").Value = "  -7.69%  "
$ws.Range("D24").Value = "'4.85"
$ws.Range("E24").Value = "  -10.09%  "
$ws.Range("D25").Value = "'5.17"
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("D26").Value = "'11.39"
$ws.Range("E26").Value = "  -7.12%  "
$ws.Range("D27").Value = "'79.15"
$ws.Range("E27").Value = "  +3.68%  "
$ws.Range("D28").Value = "'3.351.83"
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'0.154"
$ws.Range("E31").Value = "  -13.32%  "
$ws.Range("D32").Value = "'3.92"
$ws.Range("E32").Value = "  +7.99%  "
$ws.Range("D33").Value = "'8.17"
$ws.Range("E33").Value = "  -6.19%  "
$ws.Range("D34").Value = "'504.08"
$ws.Range("E34").Value = "  -10.93%  "
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "'1.85"
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "'6.72"
$ws.Range("E36").Value = "  -7.89%  "
$ws.Range("D37").Value = "'1.25"
$ws.Range("E37").Value = "  -8.79%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'22.28"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'21.85"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "'0.125"
$ws.Range("E41").Value = "  -7.69%  "
$ws.Range("D43").Value = "'1.86"
$ws.Range("E43").Value = "  -6.74%  "
$ws.Range("D44").Value = "'0.363"
$ws.Range("E44").Value = "  -8.30%  "
$ws.Range("D45").Value = "'146.82"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "'43.64"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "'166.29"
$ws.Range("E47").Value = "  -10.03%  "
$ws.Range("D48").Value = "'0.124"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").Value = "'0.721"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "'24.10"
$ws.Range("E50").Value = "  -4.73%  "
$ws.Range("D51").Value = "'1.18"
$ws.Range("E51").Value = "  -9.26%  "
